$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: push the old row8 (P8:Q8 total) and row9 (footer) down by
#        inserting 6 fresh rows right after row 7, so the new product rows
#        become 7..13 and the old rows land on 14/15 automatically.
$ws.Rows("8:13").Insert()

# --- 2. Copy row 7's formatting (borders/fills/fonts/merges-ready layout)
#        into the 6 new rows so each column keeps the same per-column style.
$ws.Range("A7:Q7").Copy()
$ws.Range("A8:Q13").PasteSpecial(-4122)
$ws.Application.CutCopyMode = 0

# --- 3. Re-create the merges for the newly inserted rows (same pattern as
#        row 7: A:B, C:G, H:K, L:M, N:O merged; P and Q stay separate).
$dataRows = 8,9,10,11,12,13
foreach ($r in $dataRows) {
    $ws.Range("A$r`:B$r").Merge()
    $ws.Range("C$r`:G$r").Merge()
    $ws.Range("H$r`:K$r").Merge()
    $ws.Range("L$r`:M$r").Merge()
    $ws.Range("N$r`:O$r").Merge()
}

# --- 4. The three text-ish columns (product name / ratio text / count text)
#        must be stored as text (numFmtId 49) instead of General so values
#        like "0:3" or "1" are not reinterpreted as numbers/dates.
$ws.Range("C7:G13").NumberFormat = "@"
$ws.Range("H7:K13").NumberFormat = "@"
$ws.Range("N7:O13").NumberFormat = "@"

# --- 5. Fill in the new item rows (A=sequence number, C=item name,
#        H=ratio text, L=count text, N=price text, P=sale-price text,
#        Q=transactions-count text).
$ws.Range("A7").Value = 1
$ws.Range("C7").Value = "AVIL 45.5MG/2ML 6 I.M. AMPS"
$ws.Range("H7").Value = "0:3"
$ws.Range("L7").Value = "1"
$ws.Range("N7").Value = "51.00"
$ws.Range("P7").Value = "51.0000"
$ws.Range("Q7").Value = "1:0"

$ws.Range("A8").Value = 2
$ws.Range("C8").Value = "CETAL COLD & FLU 20 CAPLETS"
$ws.Range("H8").Value = "0:1"
$ws.Range("L8").Value = "1"
$ws.Range("N8").Value = "36.00"
$ws.Range("P8").Value = "18.0000"
$ws.Range("Q8").Value = "0:1"

$ws.Range("A9").Value = 3
$ws.Range("C9").Value = "DECLOPHEN 75MG/3ML 3 AMPOULES"
$ws.Range("H9").Value = "3:1"
$ws.Range("L9").Value = "1"
$ws.Range("N9").Value = "36.00"
$ws.Range("P9").Value = "11.8800"
$ws.Range("Q9").Value = "0:1"

$ws.Range("A10").Value = 4
$ws.Range("C10").Value = "DEXAMETHASONE-AMRIYA 8MG/2ML 3 AMP."
$ws.Range("H10").Value = "3:2"
$ws.Range("L10").Value = "1"
$ws.Range("N10").Value = "36.00"
$ws.Range("P10").Value = "11.8800"
$ws.Range("Q10").Value = "0:1"

$ws.Range("A11").Value = 5
$ws.Range("C11").Value = "PRAVOTIN 100MG 14 SACHETS"
$ws.Range("H11").Value = "1:0"
$ws.Range("L11").Value = "1"
$ws.Range("N11").Value = "180.00"
$ws.Range("P11").Value = "180.0000"
$ws.Range("Q11").Value = "1:0"

$ws.Range("A12").Value = 6
$ws.Range("C12").Value = "RAKU-TEN EMOLLIENT CREAM 50 GM"
$ws.Range("H12").Value = "0:0"
$ws.Range("L12").Value = "1"
$ws.Range("N12").Value = "85.00"
$ws.Range("P12").Value = "85.0000"
$ws.Range("Q12").Value = "1:0"

$ws.Range("A13").Value = 7
$ws.Range("C13").Value = "سرنجات 5 سم"
$ws.Range("H13").Value = "0:0"
$ws.Range("L13").Value = "0"
$ws.Range("N13").Value = "3.00"
$ws.Range("P13").Value = "3.0000"
$ws.Range("Q13").Value = "1:0"

# --- 6. Totals row (now row 14): sum of the "sale price" column.
$ws.Range("P14").Value = 360.76

# --- 7. Footer timestamp (now row 15) refreshed to the new export time.
$ws.Range("A15").Value = "Thursday, 5 June, 2025 11:01 AM"
